$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.653636693954468
$ws.Range("B1").Value = 1.471322059631348
$ws.Range("C1").Value = 6.464095592498779
$ws.Range("D1").Value = 2.825422525405884
$ws.Range("E1").Value = 0.524304211139679
